# Update date/time field types to text on the "survey" sheet, then
# move the active selection to C7 to match the saved workbook state.

$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")

# FOL_date row: type changes from "date" to "text"
$survey.Range("C2").Value = "text"

# FOL_time_begin / FOL_time_end rows: type changes from "time" to "text"
$survey.Range("C5").Value = "text"
$survey.Range("C6").Value = "text"

# Persist the new cell selection recorded for the survey sheet
$survey.Activate()
$survey.Range("C7").Select()
